$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 109
$newRowCount = 15

for ($i = 0; $i -lt $newRowCount; $i++) {
    $destRow = $lastRow + 1 + $i
    $ws.Range("A109:O109").Copy()
    $ws.Range("A" + $destRow + ":O" + $destRow).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Row 110
$ws.Cells.Item(110, 1).Value = 44640
$ws.Cells.Item(110, 2).Value = "BriMac"
$ws.Cells.Item(110, 3).Value = 50
$ws.Cells.Item(110, 4).Value = 4
$ws.Cells.Item(110, 5).Value = "zone1"
$ws.Cells.Item(110, 6).Value = 0.5
$ws.Cells.Item(110, 7).Value = 44640.5
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 0
$ws.Cells.Item(110, 14).Value = 12
$ws.Cells.Item(110, 15).ClearContents()

# Row 111
$ws.Cells.Item(111, 1).Value = 44640
$ws.Cells.Item(111, 2).Value = "BriMac"
$ws.Cells.Item(111, 3).Value = 50
$ws.Cells.Item(111, 4).Value = 4
$ws.Cells.Item(111, 5).Value = "spur"
$ws.Cells.Item(111, 6).Value = 0.5972222222222222
$ws.Cells.Item(111, 7).Value = 44640.59722222222
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 1
$ws.Cells.Item(111, 13).Value = 0
$ws.Cells.Item(111, 14).Value = 8
$ws.Cells.Item(111, 15).ClearContents()

# Row 112
$ws.Cells.Item(112, 1).Value = 44640
$ws.Cells.Item(112, 2).Value = "BriMac"
$ws.Cells.Item(112, 3).Value = 50
$ws.Cells.Item(112, 4).Value = 4
$ws.Cells.Item(112, 5).Value = "lot3"
$ws.Cells.Item(112, 6).Value = 0.625
$ws.Cells.Item(112, 7).Value = 44640.625
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = 0
$ws.Cells.Item(112, 14).Value = 4
$ws.Cells.Item(112, 15).ClearContents()

# Row 113
$ws.Cells.Item(113, 1).Value = 44647
$ws.Cells.Item(113, 2).Value = "BriMac"
$ws.Cells.Item(113, 3).Value = 20
$ws.Cells.Item(113, 4).Value = 3
$ws.Cells.Item(113, 5).Value = "zone1"
$ws.Cells.Item(113, 6).Value = 0.5
$ws.Cells.Item(113, 7).Value = 44647.5
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(113, 15).Value = "people too numerous to count"

# Row 114
$ws.Cells.Item(114, 1).Value = 44647
$ws.Cells.Item(114, 2).Value = "BriMac"
$ws.Cells.Item(114, 3).Value = 20
$ws.Cells.Item(114, 4).Value = 3
$ws.Cells.Item(114, 5).Value = "spur"
$ws.Cells.Item(114, 6).Value = 0.5416666666666666
$ws.Cells.Item(114, 7).Value = 44647.541666666664
$ws.Cells.Item(114, 8).Value = 27
$ws.Cells.Item(114, 9).Value = 2
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 0
$ws.Cells.Item(114, 14).Value = 18
$ws.Cells.Item(114, 15).ClearContents()

# Row 115
$ws.Cells.Item(115, 1).Value = 44647
$ws.Cells.Item(115, 2).Value = "BriMac"
$ws.Cells.Item(115, 3).Value = 20
$ws.Cells.Item(115, 4).Value = 3
$ws.Cells.Item(115, 5).Value = "lot3"
$ws.Cells.Item(115, 6).Value = 0.5694444444444444
$ws.Cells.Item(115, 7).Value = 44647.569444444445
$ws.Cells.Item(115, 8).Value = 2
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 0
$ws.Cells.Item(115, 14).ClearContents()
$ws.Cells.Item(115, 15).Value = "parking lot full"

# Row 116
$ws.Cells.Item(116, 1).Value = 44654
$ws.Cells.Item(116, 2).Value = "BriMac"
$ws.Cells.Item(116, 3).Value = 100
$ws.Cells.Item(116, 4).Value = 4
$ws.Cells.Item(116, 5).Value = "zone1"
$ws.Cells.Item(116, 6).Value = 0.5
$ws.Cells.Item(116, 7).Value = 44654.5
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 0
$ws.Cells.Item(116, 14).Value = 13
$ws.Cells.Item(116, 15).ClearContents()

# Row 117
$ws.Cells.Item(117, 1).Value = 44654
$ws.Cells.Item(117, 2).Value = "BriMac"
$ws.Cells.Item(117, 3).Value = 100
$ws.Cells.Item(117, 4).Value = 4
$ws.Cells.Item(117, 5).Value = "spur"
$ws.Cells.Item(117, 6).Value = 0.6527777777777778
$ws.Cells.Item(117, 7).Value = 44654.65277777778
$ws.Cells.Item(117, 8).Value = 1
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 1
$ws.Cells.Item(117, 13).Value = 0
$ws.Cells.Item(117, 14).Value = 11
$ws.Cells.Item(117, 15).ClearContents()

# Row 118
$ws.Cells.Item(118, 1).Value = 44654
$ws.Cells.Item(118, 2).Value = "BriMac"
$ws.Cells.Item(118, 3).Value = 100
$ws.Cells.Item(118, 4).Value = 4
$ws.Cells.Item(118, 5).Value = "lot3"
$ws.Cells.Item(118, 6).Value = 0.6319444444444444
$ws.Cells.Item(118, 7).Value = 44654.631944444445
$ws.Cells.Item(118, 8).Value = 1
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = 0
$ws.Cells.Item(118, 14).Value = 5
$ws.Cells.Item(118, 15).ClearContents()

# Row 119
$ws.Cells.Item(119, 1).Value = 44675
$ws.Cells.Item(119, 2).Value = "BriMac"
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 3
$ws.Cells.Item(119, 5).Value = "zone1"
$ws.Cells.Item(119, 6).Value = 0.5833333333333334
$ws.Cells.Item(119, 7).Value = 44675.583333333336
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 13).Value = 0
$ws.Cells.Item(119, 14).Value = 17
$ws.Cells.Item(119, 15).ClearContents()

# Row 120
$ws.Cells.Item(120, 1).Value = 44675
$ws.Cells.Item(120, 2).Value = "BriMac"
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 3
$ws.Cells.Item(120, 5).Value = "spur"
$ws.Cells.Item(120, 6).Value = 0.6145833333333334
$ws.Cells.Item(120, 7).Value = 44675.614583333336
$ws.Cells.Item(120, 8).Value = 13
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = 0
$ws.Cells.Item(120, 14).Value = 8
$ws.Cells.Item(120, 15).ClearContents()

# Row 121
$ws.Cells.Item(121, 1).Value = 44675
$ws.Cells.Item(121, 2).Value = "BriMac"
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 3
$ws.Cells.Item(121, 5).Value = "lot3"
$ws.Cells.Item(121, 6).Value = 0.6354166666666666
$ws.Cells.Item(121, 7).Value = 44675.635416666664
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = 0
$ws.Cells.Item(121, 14).Value = 3
$ws.Cells.Item(121, 15).ClearContents()

# Row 122
$ws.Cells.Item(122, 1).Value = 44682
$ws.Cells.Item(122, 2).Value = "BriMac"
$ws.Cells.Item(122, 3).Value = 100
$ws.Cells.Item(122, 4).Value = 4
$ws.Cells.Item(122, 5).Value = "zone1"
$ws.Cells.Item(122, 6).Value = 0.5416666666666666
$ws.Cells.Item(122, 7).Value = 44682.541666666664
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 1
$ws.Cells.Item(122, 13).Value = 0
$ws.Cells.Item(122, 14).Value = 8
$ws.Cells.Item(122, 15).ClearContents()

# Row 123
$ws.Cells.Item(123, 1).Value = 44682
$ws.Cells.Item(123, 2).Value = "BriMac"
$ws.Cells.Item(123, 3).Value = 100
$ws.Cells.Item(123, 4).Value = 4
$ws.Cells.Item(123, 5).Value = "spur"
$ws.Cells.Item(123, 6).Value = 0.5902777777777778
$ws.Cells.Item(123, 7).Value = 44682.59027777778
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).Value = 0
$ws.Cells.Item(123, 14).Value = 5
$ws.Cells.Item(123, 15).ClearContents()

# Row 124
$ws.Cells.Item(124, 1).Value = 44682
$ws.Cells.Item(124, 2).Value = "BriMac"
$ws.Cells.Item(124, 3).Value = 100
$ws.Cells.Item(124, 4).Value = 4
$ws.Cells.Item(124, 5).Value = "lot3"
$ws.Cells.Item(124, 6).Value = 0.6180555555555556
$ws.Cells.Item(124, 7).Value = 44682.618055555555
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 13).Value = 0
$ws.Cells.Item(124, 14).Value = 0
$ws.Cells.Item(124, 15).ClearContents()

$ws.Range("M125").Select()
